$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParaByText($oldText, $newXmlInner) {
    $r = $d.Content
    $r.Find.Execute($oldText)
    if (-not $r.Find.Found) {
        throw "Text not found: $oldText"
    }
    # Expand the found range to the whole enclosing paragraph so the
    # InsertXML fragment (which includes <w:pPr> etc.) replaces cleanly.
    $para = $r.Paragraphs(1)
    $pr = $para.Range
    $xml = "<w:p $wns>$newXmlInner</w:p>"
    $pr.InsertXML($xml)
}

# 1a. Bold "meta title" paragraph near the end. The exact same sentence is
#     also the Heading1 title at the very start of the doc, so restrict the
#     search range to begin after the heading to land on the 2nd occurrence.
$tailRange = $d.Range(200, $d.Content.End)
$tailRange.Find.Execute("Play 1 Million Megaways BC for Free | Review of Prehistoric Slot")
if (-not $tailRange.Find.Found) {
    throw "Bold meta-title paragraph not found"
}
$boldXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 1 Million Megaways BC Free: A Prehistoric Adventure</w:t></w:r></w:p>'
$tailRange.Paragraphs(1).Range.InsertXML($boldXml)

# 1b. Main title heading (Heading1) - no leading empty run, plain find/replace is fine
$d.Content.Find.Execute("Play 1 Million Megaways BC for Free | Review of Prehistoric Slot", $true, $false, $false, $false, $false, $true, 1, $false, "Play 1 Million Megaways BC Free: A Prehistoric Adventure", 2)

# 2. "What we like" bullet list - reordered and reworded (preserve leading empty <w:r/>)
$listPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'

Replace-ParaByText "Impressive graphics and sound design." ("$listPPr<w:r/><w:r><w:t>Consecutive tumbling wins with multipliers up to 8x</w:t></w:r>")
Replace-ParaByText "Uses the Megaways mechanism, providing players with a huge number of ways to win." ("$listPPr<w:r/><w:r><w:t>Free Spins round with expanding multipliers and sticky wilds</w:t></w:r>")
Replace-ParaByText "Free Spins feature offers sticky wilds and expanding multipliers." ("$listPPr<w:r/><w:r><w:t>Impressive graphics and sound design</w:t></w:r>")
Replace-ParaByText "Tumbling Reels feature allows for more winning outcomes." ("$listPPr<w:r/><w:r><w:t>Huge potential for big payouts, up to 45,000x multiplied by the bet</w:t></w:r>")

# 3. "What we don't like" bullet list
Replace-ParaByText "Lower betting range limits may not appeal to some high rollers." ("$listPPr<w:r/><w:r><w:t>Limited betting range, starting at 20 cents per spin</w:t></w:r>")
Replace-ParaByText "Some players may find the prehistoric theme unappealing." ("$listPPr<w:r/><w:r><w:t>Tumbling Reels feature can become repetitive over time</w:t></w:r>")

# 4. Italic "meta description" paragraph
Replace-ParaByText "Read our review of 1 Million Megaways BC, a prehistoric-themed online slot game by Iron Dog Studio, play for free and win big with free spins." ('<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Experience big wins and impressive graphics with 1 Million Megaways BC. Play this prehistoric-themed slot for free and win big!</w:t></w:r>')
